$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "orario" (time) column F
$ws.Range("F1").Value = "orario"

# Plain h:mm (builtin format 20) cells first, so that style gets registered
# ahead of the custom h:mm;@ style used by F2
$ws.Range("F3").Value = 0.875
$ws.Range("F3").NumberFormat = "h:mm"

$ws.Range("F4").Value = 0.66666666666666663
$ws.Range("F4").NumberFormat = "h:mm"

$ws.Range("F5").Value = 0.66666666666666663
$ws.Range("F5").NumberFormat = "h:mm"

$ws.Range("F6").Value = 0.72916666666666663
$ws.Range("F6").NumberFormat = "h:mm"

$ws.Range("F7").Value = 0.79166666666666663
$ws.Range("F7").NumberFormat = "h:mm"

$ws.Range("F8").Value = 0.66666666666666663
$ws.Range("F8").NumberFormat = "h:mm"

$ws.Range("F9").Value = 0.91666666666666663
$ws.Range("F9").NumberFormat = "h:mm"

$ws.Range("F10").Value = 0.875
$ws.Range("F10").NumberFormat = "h:mm"

$ws.Range("F2").Value = 0.77083333333333337
$ws.Range("F2").NumberFormat = "h:mm;@"

# Corrected start dates for row 7 and row 10
$ws.Range("D7").Value = 44732
$ws.Range("D10").Value = 45343

# Rows 4, 5 and 7 now run for 14 days instead of 7
$ws.Range("E4").Formula = "=D4+14"
$ws.Range("E5").Formula = "=D5+14"
$ws.Range("E7").Formula = "=D7+14"

# Row 10 keeps a +7 end date, recomputed against its new start date
$ws.Range("E10").Formula = "=D10+7"

$ws.Range("F2").Select() | Out-Null
